$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.430.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.32%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'3.690.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.83%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'687.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.75%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'160.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -5.77%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'3.690.57"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -2.80%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -0.21%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("E9").Value = "'  -5.78%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("E10").Value = "'  -8.51%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("E11").Value = "'  -4.63%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("E12").Value = "'  -8.88%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("E13").Value = "'  -6.63%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'4.314.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -2.84%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'32.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -10.32%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'3.696.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -2.71%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'69.454.91"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.36%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("E18").Value = "'  -0.91%  "
$ws.Range("E18").ClearFormats()

# Row 20
$ws.Range("D20").Value = "'6.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -10.35%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'468.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -8.72%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'9.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -4.67%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("E23").Value = "'  -9.27%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'79.48"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -4.83%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'3.838.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -2.71%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").Value = "'  -10.86%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = "'10.99"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -13.44%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").Value = "'9.18"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -10.89%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("E30").Value = "'  -8.20%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("E31").Value = "'  -12.76%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("B32").Value = "'NEARProtocol"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = "'6.65"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -9.08%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("B33").Value = "'ImmutableX"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = "'2.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -10.84%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +0.13%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("D35").Value = "'26.82"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -7.93%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("E36").Value = "'  -6.45%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("E37").Value = "'  -12.02%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("E38").Value = "'  -7.88%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'2.28"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -3.49%  "
$ws.Range("E39").ClearFormats()

# Row 41
$ws.Range("E41").Value = "'  -10.19%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("E42").Value = "'  +0.05%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = "'167.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.28%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("E44").Value = "'  -6.64%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'47.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.44%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = "'2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -13.82%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("B47").Value = "'SuiNetwork"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'1.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -3.53%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("B48").Value = "'ONDO"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'1.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -4.53%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("B49").Value = "'InjectiveProtocol"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Value = "'28.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -3.28%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").Value = "'0.000279"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -8.19%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("E51").Value = "'  -9.61%  "
$ws.Range("E51").ClearFormats()

Write-Output "Applied 87 cell updates"
